$wb2 = $excel.ActiveWorkbook
$ws = $wb2.ActiveSheet

# Factor data is being trimmed from 10 data columns down to 8. Column F
# (the last column of the "Developed" block) is dropped, which folds the
# "Emerging" block (and the rest of each row) one column to the left;
# deleting entire columns keeps the surviving cells' formatting/merge
# ranges intact and auto-adjusts them instead of leaving stale styles
# behind.
$ws.Range("F1").EntireColumn.Delete()

# The new trailing column (J, originally K) is also dropped, shrinking the
# "Emerging" header's merge from F1:J1 down to F1:I1.
$ws.Range("J1").EntireColumn.Delete()

# Row 4: replace with the newly (re)computed values for the 8 remaining
# data columns.
$ws.Range("B4").Value = 0.009193250727972642
$ws.Range("C4").Value = 0.008557395566483827
$ws.Range("D4").Value = 0.005307733241375821
$ws.Range("E4").Value = 0.00768468826360108
$ws.Range("F4").Value = 0.01390734953888409
$ws.Range("G4").Value = 0.01317615080094803
$ws.Range("H4").Value = 0.01349820253513285
$ws.Range("I4").Value = 0.007444563375589609
